$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the "Foto" value for Bruno Fernandes (row 3) which was previously empty.
# The value reflects a byte-array photo being stored/serialized as its
# .NET ToString() representation: "System.Byte[]"
$ws.Range("D3").Value = "System.Byte[]"

# Resize column D (Foto) to fit the new, wider content.
$ws.Columns.Item(4).EntireColumn.AutoFit() | Out-Null
